$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "demand2" row (row 3) and the "net2" row (originally row 5,
# which becomes row 4 after the first deletion)
$ws.Rows("3").Delete()
$ws.Rows("4").Delete()

# Append the two new element groups at the bottom of the table
$ws.Range("A18").Value = "bat_with_aging1"
$ws.Range("A19").Value = "bat_with_aging2"
$ws.Range("A20").Value = "gas_boiler1"
$ws.Range("A21").Value = "gas_boiler2"
$ws.Range("B18").Value = "bat_with_aging"
$ws.Range("B19").Value = "bat_with_aging"
$ws.Range("B20").Value = "gas_boiler"
$ws.Range("B21").Value = "gas_boiler"
